$d = $word.ActiveDocument
$d.Content.Find.Execute("cartPromotor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "cart", 2)
